# Update "想去人数" (interested-people count) figures that changed
# between the two data refreshes, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 263
$ws1.Range("F4").Value = 936
$ws1.Range("F5").Value = 542

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 263
$ws4.Range("F4").Value = 936
$ws4.Range("F6").Value = 542
